$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 3500

# Row 4 updates
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

# Row 5 updates
$ws.Range("D5").Value = 44991
$ws.Range("M5").Value = 50
